# Update "想去人数" (F column) values across the four sheets to match
# the newly generated data output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 876
$ws.Cells.Item(3, 6).Value = 13875
$ws.Cells.Item(4, 6).Value = 13667
$ws.Cells.Item(6, 6).Value = 809
$ws.Cells.Item(8, 6).Value = 609
$ws.Cells.Item(9, 6).Value = 84
$ws.Cells.Item(11, 6).Value = 62
$ws.Cells.Item(13, 6).Value = 2156
$ws.Cells.Item(14, 6).Value = 118
$ws.Cells.Item(15, 6).Value = 96
$ws.Cells.Item(17, 6).Value = 132
$ws.Cells.Item(19, 6).Value = 542
$ws.Cells.Item(21, 6).Value = 424
$ws.Cells.Item(23, 6).Value = 274
$ws.Cells.Item(24, 6).Value = 845
$ws.Cells.Item(25, 6).Value = 106
$ws.Cells.Item(26, 6).Value = 8

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 141
$ws.Cells.Item(6, 6).Value = 170
$ws.Cells.Item(7, 6).Value = 1563
$ws.Cells.Item(14, 6).Value = 780
$ws.Cells.Item(15, 6).Value = 9

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 224

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 224
$ws.Cells.Item(3, 6).Value = 876
$ws.Cells.Item(4, 6).Value = 13875
$ws.Cells.Item(5, 6).Value = 13667
$ws.Cells.Item(7, 6).Value = 809
$ws.Cells.Item(9, 6).Value = 609
$ws.Cells.Item(10, 6).Value = 84
$ws.Cells.Item(12, 6).Value = 62
$ws.Cells.Item(16, 6).Value = 2156
$ws.Cells.Item(17, 6).Value = 118
$ws.Cells.Item(18, 6).Value = 96
$ws.Cells.Item(20, 6).Value = 132
$ws.Cells.Item(21, 6).Value = 141
$ws.Cells.Item(26, 6).Value = 542
$ws.Cells.Item(28, 6).Value = 424
$ws.Cells.Item(30, 6).Value = 274
$ws.Cells.Item(31, 6).Value = 845
$ws.Cells.Item(32, 6).Value = 170
$ws.Cells.Item(33, 6).Value = 1563
$ws.Cells.Item(38, 6).Value = 106
$ws.Cells.Item(41, 6).Value = 8
$ws.Cells.Item(42, 6).Value = 780
$ws.Cells.Item(43, 6).Value = 9
